$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 data: Email / Firstname / Lastname / Company -------------------
$ws.Range("A2").Value = "captain@gmail.com"
$ws.Range("B2").Value = "Archer"
$ws.Range("C2").Value = "jofra"
$ws.Range("D2").Value = "3 lions"

# Turn the new email cell into a mailto hyperlink (adds the Hyperlink cell
# style / font, xl/worksheets/_rels relationship, etc.)
$ws.Hyperlinks.Add($ws.Range("A2"), "mailto:captain@gmail.com")

# --- Column width tweaks (columns A and D get wider) -----------------------
$ws.Columns.Item(1).ColumnWidth = 25
$ws.Columns.Item(4).ColumnWidth = 19

# --- Selection moves to the newly-populated D2 cell -------------------------
$null = $ws.Range("D2").Select()
